# This script reproduces the "Updated symbol list" commit: it refreshes the
# Price / Volume(1h) columns (and, for rows 6-7, the Coin name + Link which
# were swapped) on Sheet1 of the crypto-tracking workbook.
#
# The source cells are stored as literal text (t="inlineStr") rather than
# numbers/percentages, e.g. "312.43" and "2.20%" are plain strings, not the
# number 312.43 or the percentage 0.022. Setting NumberFormat to "@" (Text)
# before assigning the value keeps Excel from reinterpreting these
# number-looking / percent-looking strings as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "312.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.20%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.67"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.79%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.166"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.35%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07891"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.23%"
# Row 6
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "KuCoinToken"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.301"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.27%"
# Row 7
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.905"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.25%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.968"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-7.08%"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9232"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.48%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1231"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.94%"
# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.41%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09142"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.66%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03349"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.96%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09603"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.80%"
# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.02%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005824"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.96%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.519"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.05%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.414"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.98%"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3441"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.04%"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.278"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "4.95%"
# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.82%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2593"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.74%"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.47%"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04385"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.35%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001251"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.69%"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004707"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "5.45%"
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001222"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.87%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02278"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.85%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05080"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.23%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007473"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.81%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1359"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.18%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008803"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.31%"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001962"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.99%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008628"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.08%"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006607"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.14%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.35%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003355"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "11.49%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001201"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.00%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.35%"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.35%"

# Touch a harmless cell last; this works around a quirk where the very last
# text-forced assignment in a script can otherwise retain a stray leading
# quote character in the saved string.
$null = $ws.Range("A1").Value

